{"js": "// Replace each three-digit x one-digit multiplication expression in the\n// document's table with its new value. The mapping is a 1:1 bijection\n// (every old expression is unique), so a simple text search-and-replace\n// per pair is unambiguous.\nconst replacements = [\n  [\"315\u00d76=1890\", \"241\u00d75=1205\"],\n  [\"133\u00d73=399\", \"533\u00d78=4264\"],\n  [\"414\u00d73=1242\", \"827\u00d73=2481\"],\n  [\"773\u00d72=1546\", \"723\u00d74=2892\"],\n  [\"951\u00d77=6657\", \"547\u00d73=1641\"],\n  [\"698\u00d79=6282\", \"523\u00d77=3661\"],\n  [\"628\u00d79=5652\", \"498\u00d77=3486\"],\n  [\"640\u00d76=3840\", \"157\u00d75=785\"],\n  [\"336\u00d76=2016\", \"660\u00d75=3300\"],\n  [\"368\u00d79=3312\", \"192\u00d74=768\"],\n  [\"327\u00d78=2616\", \"178\u00d78=1424\"],\n  [\"124\u00d78=992\", \"479\u00d79=4311\"],\n  [\"840\u00d78=6720\", \"446\u00d78=3568\"],\n  [\"774\u00d75=3870\", \"219\u00d74=876\"],\n  [\"342\u00d79=3078\", \"743\u00d74=2972\"],\n  [\"969\u00d72=1938\", \"759\u00d76=4554\"],\n  [\"596\u00d75=2980\", \"472\u00d72=944\"],\n  [\"454\u00d73=1362\", \"588\u00d78=4704\"],\n  [\"314\u00d74=1256\", \"517\u00d77=3619\"],\n  [\"937\u00d73=2811\", \"624\u00d79=5616\"],\n  [\"789\u00d78=6312\", \"712\u00d79=6408\"],\n  [\"411\u00d78=3288\", \"746\u00d78=5968\"],\n  [\"503\u00d79=4527\", \"270\u00d79=2430\"],\n  [\"644\u00d73=1932\", \"612\u00d79=5508\"],\n  [\"225\u00d74=900\", \"511\u00d78=4088\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication expression in the\n# document's table with its new value. The mapping is a 1:1 bijection\n# (every old expression is unique), so Find/Replace on each exact string\n# is unambiguous and safe to run with ReplaceAll.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"315\u00d76=1890\", \"241\u00d75=1205\"),\n    @(\"133\u00d73=399\", \"533\u00d78=4264\"),\n    @(\"414\u00d73=1242\", \"827\u00d73=2481\"),\n    @(\"773\u00d72=1546\", \"723\u00d74=2892\"),\n    @(\"951\u00d77=6657\", \"547\u00d73=1641\"),\n    @(\"698\u00d79=6282\", \"523\u00d77=3661\"),\n    @(\"628\u00d79=5652\", \"498\u00d77=3486\"),\n    @(\"640\u00d76=3840\", \"157\u00d75=785\"),\n    @(\"336\u00d76=2016\", \"660\u00d75=3300\"),\n    @(\"368\u00d79=3312\", \"192\u00d74=768\"),\n    @(\"327\u00d78=2616\", \"178\u00d78=1424\"),\n    @(\"124\u00d78=992\", \"479\u00d79=4311\"),\n    @(\"840\u00d78=6720\", \"446\u00d78=3568\"),\n    @(\"774\u00d75=3870\", \"219\u00d74=876\"),\n    @(\"342\u00d79=3078\", \"743\u00d74=2972\"),\n    @(\"969\u00d72=1938\", \"759\u00d76=4554\"),\n    @(\"596\u00d75=2980\", \"472\u00d72=944\"),\n    @(\"454\u00d73=1362\", \"588\u00d78=4704\"),\n    @(\"314\u00d74=1256\", \"517\u00d77=3619\"),\n    @(\"937\u00d73=2811\", \"624\u00d79=5616\"),\n    @(\"789\u00d78=6312\", \"712\u00d79=6408\"),\n    @(\"411\u00d78=3288\", \"746\u00d78=5968\"),\n    @(\"503\u00d79=4527\", \"270\u00d79=2430\"),\n    @(\"644\u00d73=1932\", \"612\u00d79=5508\"),\n    @(\"225\u00d74=900\", \"511\u00d78=4088\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    # wdReplaceAll = 2 (every exact expression is unique so this is a\n    # single, unambiguous hit per pair, matching ReplaceOne semantics too).\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
